$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.05
$ws.Range("I2").Value = 2.92

$ws.Range("J3").Value = 3.6
$ws.Range("Q3").Value = 2.02
$ws.Range("R3").Value = 1.36
$ws.Range("S3").Value = 3.4

$ws.Range("G4").Value = 2.58
$ws.Range("H4").Value = 2.92
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 3.65
$ws.Range("K4").Value = 4.5
$ws.Range("P4").Value = 2.22
$ws.Range("Q4").Value = 1.68

$ws.Range("H5").Value = 1.76
$ws.Range("I5").Value = 1.87
$ws.Range("P5").Value = 2.82

$ws.Range("H6").Value = 2.48

$wb.Save()
